# Update "想去人数" (people interested) counts on the sheets that carry
# the exhibition data: "展览" and "全部类型" (sheet1 & sheet4).
#   F2: 283 -> 287
#   F4: 27  -> 31
#   F5: 265 -> 266

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 287
    $ws.Range("F4").Value = 31
    $ws.Range("F5").Value = 266
}
